$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Data table edits ---
# Row 14
$ws.Range("N14").Value = -100

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -38.461538461538
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 14.285714285714
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = -11.111111111111
$ws.Range("N16").Value = -77.777777777777

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -46.666666666666
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = -25
$ws.Range("L17").Value = -14.285714285714
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -60

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 3
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = -62.5
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -92.5

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 62.5
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 11
$ws.Range("K19").Value = 27.272727272727
$ws.Range("L19").Value = -6.666666666666
$ws.Range("M19").Value = 7.692307692307
$ws.Range("N19").Value = -64.102564102564

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 1
$ws.Range("L20").Value = -87.5
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -96.153846153846

# Row 21
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 49
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = -3.92156862745
$ws.Range("I21").Value = 34
$ws.Range("J21").Value = 27
$ws.Range("K21").Value = 25.925925925925
$ws.Range("L21").Value = -22.727272727272
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -78.75

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("L22").Value = 100

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = 20
$ws.Range("L23").Value = -14.285714285714
$ws.Range("M23").Value = 20

# Row 24
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 22
$ws.Range("G24").Value = 28
$ws.Range("H24").Value = -21.428571428571
$ws.Range("I24").Value = 14
$ws.Range("J24").Value = 16
$ws.Range("K24").Value = -12.5
$ws.Range("L24").Value = -26.315789473684
$ws.Range("M24").Value = -46.153846153846

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 4
$ws.Range("K25").Value = -25
$ws.Range("L25").Value = -25

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 14.285714285714
$ws.Range("I26").Value = 19
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = 26.666666666666
$ws.Range("L26").Value = 90
$ws.Range("M26").Value = 137.5

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 1
$ws.Range("L28").Value = -66.666666666666

# Row 29
$ws.Range("L29").Value = -100
$ws.Range("N29").Value = -100

# Row 30
$ws.Range("L30").Value = -100
$ws.Range("N30").Value = -100

# --- Format fix-ups (restore correct style/number-format after type changes) ---
$ws.Range("H14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("L29").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("N29").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("N30").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
